# Auto-generated Excel COM-interop script applying the Ultros_Profits.xlsx commit diff
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 9165.5
$ws.Range("J40").Value = 9165.5
$ws.Range("L40").Value = 9165.5
$ws.Range("N40").Value = -9515.5
$ws.Range("H43").Value = 4611.125
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7931
$ws.Range("H53").Value = 1178.5
$ws.Range("J53").Value = 1238
$ws.Range("L53").Value = 1238
$ws.Range("N53").Value = -2512
$ws.Range("H132").Value = 12682.866
$ws.Range("I132").Value = 1179.7441
$ws.Range("J132").Value = 260000
$ws.Range("K132").Value = 3539.2323
$ws.Range("L132").Value = 780000
$ws.Range("M132").Value = -1009.2323
$ws.Range("N132").Value = -785060
# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 44438.625
$ws.Range("J2").Value = 6499.3335
$ws.Range("L2").Value = 6499.3335
$ws.Range("N2").Value = -6725.3335
$ws.Range("H32").Value = 2198.8333
$ws.Range("I32").Value = 2280.976
$ws.Range("K32").Value = 2280.976
$ws.Range("M32").Value = -1993.976
$ws.Range("H61").Value = 20848.75
$ws.Range("I61").Value = 7300
$ws.Range("K61").Value = 7300
$ws.Range("M61").Value = -7088
$ws.Range("H74").Value = 9486.166999999999
$ws.Range("I74").Value = 6974.3335
$ws.Range("K74").Value = 6974.3335
$ws.Range("M74").Value = -6100.3335
$ws.Range("H77").Value = 9486.166999999999
$ws.Range("I77").Value = 6974.3335
$ws.Range("K77").Value = 34871.6675
$ws.Range("M77").Value = -30503.6675
$ws.Range("H97").Value = 4437.778
$ws.Range("I97").Value = 2570
$ws.Range("J97").Value = 10975
$ws.Range("K97").Value = 2570
$ws.Range("L97").Value = 10975
$ws.Range("M97").Value = -2074
$ws.Range("N97").Value = -11967
$ws.Range("H110").Value = 5862.1665
$ws.Range("I110").Value = 5770.3335
$ws.Range("J110").Value = 6137.6665
$ws.Range("K110").Value = 5770.3335
$ws.Range("L110").Value = 6137.6665
$ws.Range("M110").Value = -3725.3335
$ws.Range("N110").Value = -10227.6665
$ws.Range("H116").Value = 44438.625
$ws.Range("J116").Value = 6499.3335
$ws.Range("L116").Value = 6499.3335
$ws.Range("N116").Value = -11087.3335
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
$ws.Range("H136").Value = 20848.75
$ws.Range("I136").Value = 7300
$ws.Range("K136").Value = 21900
$ws.Range("M136").Value = -19350
# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 44438.625
$ws.Range("J3").Value = 6499.3335
$ws.Range("L3").Value = 6499.3335
$ws.Range("N3").Value = -6727.3335
$ws.Range("H40").Value = 45831
$ws.Range("J40").Value = 45831
$ws.Range("L40").Value = 45831
$ws.Range("N40").Value = -46361
$ws.Range("H64").Value = 719
$ws.Range("J64").Value = 612.1111
$ws.Range("L64").Value = 612.1111
$ws.Range("N64").Value = -1062.1111
$ws.Range("H67").Value = 719
$ws.Range("J67").Value = 612.1111
$ws.Range("L67").Value = 612.1111
$ws.Range("N67").Value = -2172.1111
$ws.Range("H96").Value = 21332.666
$ws.Range("I96").Value = 5799.4
$ws.Range("J96").Value = 98999
$ws.Range("K96").Value = 5799.4
$ws.Range("L96").Value = 98999
$ws.Range("M96").Value = -3053.4
$ws.Range("N96").Value = -104491
$ws.Range("H105").Value = 1763.6666
$ws.Range("I105").Value = 990.5625
$ws.Range("J105").Value = 3309.875
$ws.Range("K105").Value = 990.5625
$ws.Range("L105").Value = 3309.875
$ws.Range("M105").Value = 756.4375
$ws.Range("N105").Value = -6803.875
$ws.Range("H134").Value = 2575.8
$ws.Range("I134").Value = 2593
$ws.Range("J134").Value = 2507
$ws.Range("K134").Value = 7779
$ws.Range("L134").Value = 7521
$ws.Range("M134").Value = -5244
$ws.Range("N134").Value = -12591
# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 5557
$ws.Range("I16").Value = 3179.8
$ws.Range("J16").Value = 11500
$ws.Range("K16").Value = 3179.8
$ws.Range("L16").Value = 11500
$ws.Range("M16").Value = -2892.8
$ws.Range("N16").Value = -12074
$ws.Range("H31").Value = 2942.6086
$ws.Range("I31").Value = 1800.6154
$ws.Range("K31").Value = 1800.6154
$ws.Range("M31").Value = -1505.6154
$ws.Range("H34").Value = 2942.6086
$ws.Range("I34").Value = 1800.6154
$ws.Range("K34").Value = 1800.6154
$ws.Range("M34").Value = -1598.6154
$ws.Range("H107").Value = 5003.364
$ws.Range("I107").Value = 207.58333
$ws.Range("J107").Value = 10758.3
$ws.Range("K107").Value = 207.58333
$ws.Range("L107").Value = 10758.3
$ws.Range("M107").Value = 1712.41667
$ws.Range("N107").Value = -14598.3
$ws.Range("H113").Value = 5557
$ws.Range("I113").Value = 3179.8
$ws.Range("J113").Value = 11500
$ws.Range("K113").Value = 3179.8
$ws.Range("L113").Value = 11500
$ws.Range("M113").Value = -1009.8
$ws.Range("N113").Value = -15840
$ws.Range("H132").Value = 2469.7222
$ws.Range("I132").Value = 2117.3333
$ws.Range("K132").Value = 6351.999899999999
$ws.Range("M132").Value = -3821.999899999999
$ws.Range("H134").Value = 9112.286
$ws.Range("I134").Value = 9112.286
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 27336.858
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -24801.858
$ws.Range("N134").ClearContents()
# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H3").Value = 5599.7646
$ws.Range("I3").Value = 4699.75
$ws.Range("K3").Value = 14099.25
$ws.Range("M3").Value = -13987.25
$ws.Range("H59").Value = 33920.168
$ws.Range("I59").Value = 673.25
$ws.Range("J59").Value = 50543.625
$ws.Range("K59").Value = 2019.75
$ws.Range("L59").Value = 151630.875
$ws.Range("M59").Value = -1479.75
$ws.Range("N59").Value = -152710.875
$ws.Range("H68").Value = 2100.3333
$ws.Range("J68").Value = 2100.3333
$ws.Range("L68").Value = 6300.999899999999
$ws.Range("N68").Value = -7922.999899999999
$ws.Range("H71").Value = 2100.3333
$ws.Range("J71").Value = 2100.3333
$ws.Range("L71").Value = 18902.9997
$ws.Range("N71").Value = -27014.9997
$ws.Range("H81").Value = 3000
$ws.Range("I81").Value = 3000
$ws.Range("K81").Value = 9000
$ws.Range("M81").Value = -7877
$ws.Range("H84").Value = 3000
$ws.Range("I84").Value = 3000
$ws.Range("K84").Value = 27000
$ws.Range("M84").Value = -21384
$ws.Range("H122").Value = 5288.6
$ws.Range("J122").Value = 6482.5
$ws.Range("L122").Value = 58342.5
$ws.Range("N122").Value = -63242.5
$ws.Range("H133").Value = 8343
$ws.Range("I133").Value = 7514.5
$ws.Range("K133").Value = 22543.5
$ws.Range("M133").Value = -17483.5
$ws.Range("H140").Value = 626772.6
$ws.Range("I140").Value = 668324.2
$ws.Range("J140").Value = 3499
$ws.Range("K140").Value = 2004972.6
$ws.Range("L140").Value = 10497
$ws.Range("M140").Value = -1999792.6
$ws.Range("N140").Value = -20857
# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 4166826
$ws.Range("I2").Value = 60.166668
$ws.Range("J2").Value = 8333592
$ws.Range("K2").Value = 60.166668
$ws.Range("L2").Value = 8333592
$ws.Range("M2").Value = 52.833332
$ws.Range("N2").Value = -8333818
$ws.Range("H27").Value = 4500
$ws.Range("J27").Value = 4500
$ws.Range("L27").Value = 4500
$ws.Range("N27").Value = -4832
$ws.Range("H70").Value = 107840.63
$ws.Range("I70").Value = 165771.42
$ws.Range("K70").Value = 165771.42
$ws.Range("M70").Value = -165501.42
$ws.Range("H73").Value = 107840.63
$ws.Range("I73").Value = 165771.42
$ws.Range("K73").Value = 165771.42
$ws.Range("M73").Value = -164835.42
$ws.Range("H132").Value = 9655.77
$ws.Range("I132").Value = 9264.764999999999
$ws.Range("K132").Value = 27794.295
$ws.Range("M132").Value = -25264.295
# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H80").Value = 75000
$ws.Range("J80").Value = 75000
$ws.Range("L80").Value = 75000
$ws.Range("N80").Value = -77246
$ws.Range("H83").Value = 75000
$ws.Range("J83").Value = 75000
$ws.Range("L83").Value = 225000
$ws.Range("N83").Value = -236232
$ws.Range("H132").Value = 3851
$ws.Range("I132").Value = 2704
$ws.Range("K132").Value = 8112
$ws.Range("M132").Value = -5582
# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H44").Value = 36829.332
$ws.Range("I44").Value = 22993
$ws.Range("J44").Value = 43747.5
$ws.Range("K44").Value = 22993
$ws.Range("L44").Value = 43747.5
$ws.Range("M44").Value = -22439
$ws.Range("N44").Value = -44855.5
$ws.Range("H96").Value = 125625.75
$ws.Range("I96").Value = 125625.75
$ws.Range("K96").Value = 125625.75
$ws.Range("M96").Value = -124252.75
$ws.Range("H113").Value = 1098.0714
$ws.Range("I113").Value = 635.2857
$ws.Range("J113").Value = 1560.8572
$ws.Range("K113").Value = 1905.8571
$ws.Range("L113").Value = 4682.571599999999
$ws.Range("M113").Value = 264.1428999999998
$ws.Range("N113").Value = -9022.571599999999
$ws.Range("H122").Value = 1237.75
$ws.Range("I122").Value = 937.4737
$ws.Range("J122").Value = 2378.8
$ws.Range("K122").Value = 2812.4211
$ws.Range("L122").Value = 7136.400000000001
$ws.Range("M122").Value = -362.4211
$ws.Range("N122").Value = -12036.4
$ws.Range("H136").Value = 123362.375
$ws.Range("I136").Value = 161983.33
$ws.Range("J136").Value = 7499.5
$ws.Range("K136").Value = 485949.99
$ws.Range("L136").Value = 22498.5
$ws.Range("M136").Value = -483399.99
$ws.Range("N136").Value = -27598.5
